# Add 2022-Q1 data
# 1. Insert a new worksheet "2022-Q1" positioned right before the "总计" sheet,
#    using the "2021-Q4" sheet as a formatting template (same header/row style).
# 2. Insert a new first-data-row into the "总计" sheet summarizing the new quarter,
#    shifting the existing rows down and re-indexing the helper index column.

$wb = $excel.ActiveWorkbook

$templateWs = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" worksheet right before "总计" ---
# NOTE: worksheet handles in this runtime are index-bound, not object
# identities. Inserting a sheet shifts the index of everything after it, so
# "总计" must be re-looked-up (by name) *after* the insert rather than reused
# from a variable captured beforehand.
$totalWsBeforeInsert = $wb.Worksheets.Item("总计")
$newWs = $wb.Worksheets.Add($totalWsBeforeInsert)
$newWs.Name = "2022-Q1"

# Copy the cell formatting (fonts/borders/alignment) from the template sheet so
# the new sheet visually matches the other per-quarter sheets.
$templateWs.Range("B1:H4").Copy()
$newWs.Range("B1:H4").PasteSpecial(-4122)
$templateWs.Range("A2:A4").Copy()
$newWs.Range("A2:A4").PasteSpecial(-4122)

# Header row
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Make sure the fund-code / figure columns are stored as text, matching the
# source data (values such as "002810" must keep their leading zero).
$newWs.Range("B2:G4").NumberFormat = "@"

# Row 2
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "233009"
$newWs.Range("C2").Value = "大摩多因子精选策略混合"
$newWs.Range("D2").Value = "6.77"
$newWs.Range("E2").Value = "89.73"
$newWs.Range("F2").Value = "1.07"
$newWs.Range("G2").Value = "0.0724"
$newWs.Range("H2").Value = 8

# Row 3
$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "002810"
$newWs.Range("C3").Value = "金信转型创新成长灵活配置混合"
$newWs.Range("D3").Value = "0.18"
$newWs.Range("E3").Value = "81.12"
$newWs.Range("F3").Value = "4.22"
$newWs.Range("G3").Value = "0.0076"
$newWs.Range("H3").Value = 4

# Row 4
$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "004536"
$newWs.Range("C4").Value = "嘉实中小企业量化活力灵活配置混合"
$newWs.Range("D4").Value = "0.17"
$newWs.Range("E4").Value = "90.06"
$newWs.Range("F4").Value = "1.39"
$newWs.Range("G4").Value = "0.0024"
$newWs.Range("H4").Value = 10

# --- 2. Update the "总计" (totals) sheet with the new quarter ---
# Re-fetch "总计" now that it has shifted from index 4 to index 5.
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()
$totalWs.Range("B2:D2").ClearFormats()

# Give the new A2 index cell the same style used by the other index cells.
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 3
$totalWs.Range("D2").Value = 0.08

# Re-number the index column for the rows that shifted down.
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
